# daily auto push: 2026-01-15 02:26 UTC
#
# A new observation (2026/01/15, 木, hour 9, rank 201) was recorded between
# the existing 2026/01/15 06:00 row and the 2026/12/29 run, so insert one
# new row at row 654 and push everything from the old row 654 down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 654 - this shifts the old rows 654:695
# down to 655:696 and grows the sheet dimension to A1:D696 automatically.
$ws.Rows.Item(654).Insert()

# Column A holds plain text dates (e.g. "2026/01/15"), not real date
# serials. Assigning a date-shaped string straight to .Value would let
# Excel auto-convert it into a date value, so force the cell to Text
# first, write the literal string, then drop the explicit number format
# again so the cell matches the plain (unstyled) cells around it.
$ws.Range("A654").NumberFormat = "@"
$ws.Range("A654").Value = "2026/01/15"
$ws.Range("A654").ClearFormats()

$ws.Range("B654").Value = "木"
$ws.Range("C654").Value = 9
$ws.Range("D654").Value = 201
